# Update "想去人数" (want-to-go count) values in column F
# for sheets "展览" and "全部类型", as per the commit's data refresh.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 135
    3  = 50
    5  = 89
    7  = 1251
    8  = 1529
    9  = 338
    12 = 146
    14 = 62
    15 = 106
    17 = 299
    19 = 1723
    23 = 665
    25 = 332
    26 = 4152
    28 = 265
    29 = 1082
    30 = 483
    32 = 518
    34 = 236
    37 = 13
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
